# Refresh the COVID-19 "paises" sheet with the 23:26 snapshot of 19 June 2020.
# The table (rows 4-219) is sorted descending by "Casos totales" (col B), so
# a handful of countries whose totals grew past a neighbour's now occupy a
# different row than before; those rows get both a new country name (col A)
# and new stats (cols B-H). Every other changed row keeps its country and
# just gets refreshed stats.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp title in A1
$ws.Range("A1").Value = "Datos actualizados a 19 de Junio de 2020 a las 23:26"

# Row 4
$ws.Range("B4").Value = 2291145
$ws.Range("C4").Value = 27494
$ws.Range("D4").Value = 935478
$ws.Range("E4").Value = 1234359
$ws.Range("G4").Value = 620
$ws.Range("H4").Value = 121308

# Row 5
$ws.Range("B5").Value = 1032913
$ws.Range("C5").Value = 49554
$ws.Range("E5").Value = 463599
$ws.Range("G5").Value = 1085
$ws.Range("H5").Value = 48954

# Row 9
$ws.Range("H9").Value = 28315

# Row 10
$ws.Range("B10").Value = 247925
$ws.Range("C10").Value = 3537
$ws.Range("D10").Value = 135520
$ws.Range("E10").Value = 104745
$ws.Range("G10").Value = 199
$ws.Range("H10").Value = 7660

# Row 11
$ws.Range("D11").Value = 181907
$ws.Range("E11").Value = 21543

# Row 12
$ws.Range("D12").Value = 191491
$ws.Range("E12").Value = 35809

# Row 14
$ws.Range("B14").Value = 190660
$ws.Range("C14").Value = 534
$ws.Range("E14").Value = 7300
$ws.Range("G14").Value = 14
$ws.Range("H14").Value = 8960

# Row 50
$ws.Range("B50").Value = 20916
$ws.Range("C50").Value = 486
$ws.Range("D50").Value = 15287
$ws.Range("E50").Value = 5572

# Row 81
$ws.Range("B81").Value = 4904
$ws.Range("C81").Value = 63
$ws.Range("D81").Value = 3522
$ws.Range("E81").Value = 1355
$ws.Range("G81").Value = 1
$ws.Range("H81").Value = 27

# Rows 115-118: Guinea-Bisau's total overtook Eslovaquia's, so it now sits
# right after Eslovaquia (row 115), pushing Eslovenia/Libano/Nueva Zelanda
# each down one row (116-118).
# Row 115
$ws.Range("A115").Value = "Guinea-Bisau"
$ws.Range("B115").Value = 1541
$ws.Range("C115").Value = 49
$ws.Range("D115").Value = 153
$ws.Range("E115").Value = 1371
$ws.Range("G115").Value = 2
$ws.Range("H115").Value = 17

# Row 116
$ws.Range("A116").Value = "Eslovenia"
$ws.Range("B116").Value = 1513
$ws.Range("C116").Value = 2
$ws.Range("D116").Value = 1359
$ws.Range("E116").Value = 45
$ws.Range("H116").Value = 109

# Row 117
$ws.Range("A117").Value = "Libano"
$ws.Range("B117").Value = 1510
$ws.Range("C117").Value = 15
$ws.Range("D117").Value = 960
$ws.Range("E117").Value = 518
$ws.Range("H117").Value = 32

# Row 118
$ws.Range("A118").Value = "Nueva Zelanda"
$ws.Range("B118").Value = 1507
$ws.Range("D118").Value = 1482
$ws.Range("E118").Value = 3
$ws.Range("H118").Value = 22

# Row 148
$ws.Range("B148").Value = 620
$ws.Range("C148").Value = 28
$ws.Range("D148").Value = 91
$ws.Range("E148").Value = 521

# Rows 156-157: Montenegro overtook Vietnam.
# Row 156
$ws.Range("A156").Value = "Montenegro"
$ws.Range("B156").Value = 355
$ws.Range("C156").Value = 18
$ws.Range("D156").Value = 315
$ws.Range("E156").Value = 31
$ws.Range("H156").Value = 9

# Row 157
$ws.Range("A157").Value = "Vietnam"
$ws.Range("B157").Value = 349
$ws.Range("C157").Value = 7
$ws.Range("D157").Value = 326
$ws.Range("E157").Value = 23
$ws.Range("H157").Value = 0

# Row 162
$ws.Range("B162").Value = 236
$ws.Range("C162").Value = 15
$ws.Range("E162").Value = 124

# Rows 202-203: Dominica overtook Fiyi (same totals, order swaps).
# Row 202
$ws.Range("A202").Value = "Dominica"

# Row 203
$ws.Range("A203").Value = "Fiyi"

# Rows 208-209: Islas Turcas y Caicos overtook Santa Sede.
# Row 208
$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("D208").Value = 11
$ws.Range("H208").Value = 1

# Row 209
$ws.Range("A209").Value = "Santa Sede"
$ws.Range("D209").Value = 12
$ws.Range("H209").Value = 0
